# Update the DataPruebas sheet test-case row (row 3) with the real
# "Password Corta" test case data, replacing the old placeholder values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataPruebas")

$ws.Range("A3").Value = "CP002_Registro_Fallido_Password_Corta"
$ws.Range("D3").Value = " "
$ws.Range("E3").Value = " "
$ws.Range("F3").Value = " "
$ws.Range("B3").Value = "  "
$ws.Range("C3").Value = "  "
$ws.Range("H3").Value = "  "
$ws.Range("G3").Value = "Junio"
$ws.Range("I3").Value = "Tu contraseña es demasiado corta."
$ws.Range("J3").ClearContents()

# Widen column I slightly to fit the new (longer) text.
$ws.Columns.Item(9).ColumnWidth = 29

# Move the active selection to D19, matching where the author left the cursor.
$ws.Range("D19").Select()
